$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.042315208657007
$ws.Range("D2").Value = 1.036800625962593
$ws.Range("E2").Value = 1.04940325582198
$ws.Range("F2").Value = 1.057760421188101
$ws.Range("I2").Value = 1.039793709644768
$ws.Range("J2").Value = 1.047391692388213
$ws.Range("K2").Value = 1.039593331127416
$ws.Range("L2").Value = 1.052160411215785
$ws.Range("M2").Value = 1.060494508644685
$ws.Range("N2").Value = 1.019621067833031
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.04350435527642
$ws.Range("D3").Value = 1.037363547642832
$ws.Range("E3").Value = 1.050497898448737
$ws.Range("F3").Value = 1.059027304634823
$ws.Range("I3").Value = 1.040061590554108
$ws.Range("J3").Value = 1.048226035729607
$ws.Range("K3").Value = 1.039966583463774
$ws.Range("L3").Value = 1.053066563696003
$ws.Range("M3").Value = 1.061574142107846
$ws.Range("N3").Value = 1.019903609287751
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.044273430608172
$ws.Range("D4").Value = 1.037727525350693
$ws.Range("E4").Value = 1.051206245785404
$ws.Range("F4").Value = 1.059847393168399
$ws.Range("I4").Value = 1.040233309818786
$ws.Range("J4").Value = 1.04876500980023
$ws.Range("K4").Value = 1.040207129149878
$ws.Range("L4").Value = 1.053652354057521
$ws.Range("M4").Value = 1.062272509390017
$ws.Range("N4").Value = 1.02008598581088
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.044596660283887
$ws.Range("D5").Value = 1.037880476017008
$ws.Range("E5").Value = 1.051504045809342
$ws.Range("F5").Value = 1.060192239027051
$ws.Range("I5").Value = 1.040305113380398
$ws.Range("J5").Value = 1.048991379598007
$ws.Range("K5").Value = 1.040308021119713
$ws.Range("L5").Value = 1.053898489504403
$ws.Range("M5").Value = 1.062566049809865
$ws.Range("N5").Value = 1.020162550205713
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.04465092676918
$ws.Range("D6").Value = 1.037906153237088
$ws.Range("E6").Value = 1.051554048410385
$ws.Range("F6").Value = 1.060250144938766
$ws.Range("I6").Value = 1.040317146808789
$ws.Range("J6").Value = 1.049029375510372
$ws.Range("K6").Value = 1.040324947625668
$ws.Range("L6").Value = 1.053939809088383
$ws.Range("M6").Value = 1.0626153334581
$ws.Range("N6").Value = 1.02017539945487
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04427774996863
$ws.Range("D7").Value = 1.0377295693437
$ws.Range("E7").Value = 1.05121022496003
$ws.Range("F7").Value = 1.059852000696811
$ws.Range("I7").Value = 1.040234270782951
$ws.Range("J7").Value = 1.048768035407718
$ws.Range("K7").Value = 1.040208478191043
$ws.Range("L7").Value = 1.05365564344413
$ws.Range("M7").Value = 1.062276431898789
$ws.Range("N7").Value = 1.020087009287011
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.04271716655607
$ws.Range("D8").Value = 1.036990923424569
$ws.Range("E8").Value = 1.049773187754193
$ws.Range("F8").Value = 1.05818850295723
$ws.Range("I8").Value = 1.039884576210931
$ws.Range("J8").Value = 1.047673849970456
$ws.Range("K8").Value = 1.039719674858042
$ws.Range("L8").Value = 1.052466764394224
$ws.Range("M8").Value = 1.06085942379625
$ws.Range("N8").Value = 1.019716646677932
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.039964209072073
$ws.Range("D9").Value = 1.035687302469487
$ws.Range("E9").Value = 1.047241202055711
$ws.Range("F9").Value = 1.055259665297458
$ws.Range("I9").Value = 1.039255978611594
$ws.Range("J9").Value = 1.045738802657128
$ws.Range("K9").Value = 1.038850900644257
$ws.Range("L9").Value = 1.050367540426369
$ws.Range("M9").Value = 1.058360669571199
$ws.Range("N9").Value = 1.01906058818635
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.038126760129324
$ws.Range("D10").Value = 1.034816904084407
$ws.Range("E10").Value = 1.04555330970536
$ws.Range("F10").Value = 1.053308657039039
$ws.Range("I10").Value = 1.038828577559311
$ws.Range("J10").Value = 1.044444030248121
$ws.Range("K10").Value = 1.038266735053468
$ws.Range("L10").Value = 1.04896511911305
$ws.Range("M10").Value = 1.056693532682821
$ws.Range("N10").Value = 1.01862089161656
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.037330588411209
$ws.Range("D11").Value = 1.034439709186721
$ws.Range("E11").Value = 1.044822439405049
$ws.Range("F11").Value = 1.052464192978178
$ws.Range("I11").Value = 1.038641529329598
$ws.Range("J11").Value = 1.043882239517319
$ws.Range("K11").Value = 1.038012606653993
$ws.Range("L11").Value = 1.048357141916208
$ws.Range("M11").Value = 1.055971315263222
$ws.Range("N11").Value = 1.018429942855686
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.037034770551314
$ws.Range("D12").Value = 1.03429955686657
$ws.Range("E12").Value = 1.044550959753744
$ws.Range("F12").Value = 1.05215056967872
$ws.Range("I12").Value = 1.038571753479902
$ws.Range("J12").Value = 1.043673391855925
$ws.Range("K12").Value = 1.037918034864038
$ws.Range("L12").Value = 1.048131202517521
$ws.Range("M12").Value = 1.055702999212001
$ws.Range("N12").Value = 1.01835893187749
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.037098228338032
$ws.Range("D13").Value = 1.034329622055139
$ws.Range("E13").Value = 1.044609193171557
$ws.Range("F13").Value = 1.052217840791925
$ws.Range("I13").Value = 1.038586734116627
$ws.Range("J13").Value = 1.043718198265808
$ws.Range("K13").Value = 1.037938328848639
$ws.Range("L13").Value = 1.048179672269169
$ws.Range("M13").Value = 1.055760556314189
$ws.Range("N13").Value = 1.018374167785597
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.03730613774732
$ws.Range("D14").Value = 1.034428125072964
$ws.Range("E14").Value = 1.044799998860086
$ws.Range("F14").Value = 1.052438267806453
$ws.Range("I14").Value = 1.038635767712611
$ws.Range("J14").Value = 1.043864979653592
$ws.Range("K14").Value = 1.038004792933648
$ws.Range("L14").Value = 1.04833846793948
$ws.Range("M14").Value = 1.055949137253228
$ws.Range("N14").Value = 1.018424074784068
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.03743422636026
$ws.Range("D15").Value = 1.034488810046861
$ws.Range("E15").Value = 1.044917560244702
$ws.Range("F15").Value = 1.052574086495019
$ws.Range("I15").Value = 1.038665939460101
$ws.Range("J15").Value = 1.043955393468165
$ws.Range("K15").Value = 1.038045720177465
$ws.Range("L15").Value = 1.048436292625185
$ws.Range("M15").Value = 1.056065321202283
$ws.Range("N15").Value = 1.018454812976435
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.038179587747436
$ws.Range("D16").Value = 1.034841930869288
$ws.Range("E16").Value = 1.045601814985224
$ws.Range("F16").Value = 1.053364708221711
$ws.Range("I16").Value = 1.03884094958315
$ws.Range("J16").Value = 1.044481290186952
$ws.Range("K16").Value = 1.038283575858406
$ws.Range("L16").Value = 1.049005453270176
$ws.Range("M16").Value = 1.056741456635065
$ws.Range("N16").Value = 1.018633552483496
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.038646985810553
$ws.Range("D17").Value = 1.035063352650953
$ws.Range("E17").Value = 1.046031028294044
$ws.Range("F17").Value = 1.053860732629976
$ws.Range("I17").Value = 1.038950198349312
$ws.Range("J17").Value = 1.044810863546665
$ws.Range("K17").Value = 1.038432460455647
$ws.Range("L17").Value = 1.049362279014987
$ws.Range("M17").Value = 1.057165487233743
$ws.Range("N17").Value = 1.018745521591502
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.038919558854289
$ws.Range("D18").Value = 1.035192474594813
$ws.Range("E18").Value = 1.046281381094754
$ws.Range("F18").Value = 1.054150087781443
$ws.Range("I18").Value = 1.039013730208177
$ws.Range("J18").Value = 1.045002987663223
$ws.Range("K18").Value = 1.038519188374241
$ws.Range("L18").Value = 1.049570340114941
$ws.Range("M18").Value = 1.057412784800769
$ws.Range("N18").Value = 1.018810777538475
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.039012490414978
$ws.Range("D19").Value = 1.035236496791529
$ws.Range("E19").Value = 1.046366744968029
$ws.Range("F19").Value = 1.054248755988095
$ws.Range("I19").Value = 1.039035360523133
$ws.Range("J19").Value = 1.045068478331054
$ws.Range("K19").Value = 1.03854874101857
$ws.Range("L19").Value = 1.049641271860196
$ws.Range("M19").Value = 1.057497101455705
$ws.Range("N19").Value = 1.018833019026403
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.038596843854517
$ws.Range("D20").Value = 1.035039599249745
$ws.Range("E20").Value = 1.04598497778699
$ws.Range("F20").Value = 1.053807510559644
$ws.Range("I20").Value = 1.038938496758784
$ws.Range("J20").Value = 1.04477551487966
$ws.Range("K20").Value = 1.038416498332344
$ws.Range("L20").Value = 1.049324002171369
$ws.Range("M20").Value = 1.057119996119541
$ws.Range("N20").Value = 1.018733513927487
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.037244915953478
$ws.Range("D21").Value = 1.034399119635145
$ws.Range("E21").Value = 1.044743811405406
$ws.Range("F21").Value = 1.052373356229861
$ws.Range("I21").Value = 1.038621336763519
$ws.Range("J21").Value = 1.043821760996281
$ws.Range("K21").Value = 1.037985225806684
$ws.Range("L21").Value = 1.048291709609952
$ws.Range("M21").Value = 1.055893606326746
$ws.Range("N21").Value = 1.018409380741785
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.036394416640135
$ws.Range("D22").Value = 1.033996162708268
$ws.Range("E22").Value = 1.043963429159635
$ws.Range("F22").Value = 1.05147192459246
$ws.Range("I22").Value = 1.038420202542696
$ws.Range("J22").Value = 1.043221093135096
$ws.Range("K22").Value = 1.037713042605717
$ws.Range("L22").Value = 1.047642031460221
$ws.Range("M22").Value = 1.055122222620098
$ws.Range("N22").Value = 1.018205098535913
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.036845329558809
$ws.Range("D23").Value = 1.034209802390307
$ws.Range("E23").Value = 1.044377126201429
$ws.Range("F23").Value = 1.051949764912849
$ws.Range("I23").Value = 1.038526991008533
$ws.Range("J23").Value = 1.043539614271855
$ws.Range("K23").Value = 1.037857429161127
$ws.Range("L23").Value = 1.0479864988181
$ws.Range("M23").Value = 1.055531177129241
$ws.Range("N23").Value = 1.018313438672
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.038619501007847
$ws.Range("D24").Value = 1.035050332481095
$ws.Range("E24").Value = 1.04600578602918
$ws.Range("F24").Value = 1.053831559221933
$ws.Range("I24").Value = 1.038943784794389
$ws.Range("J24").Value = 1.044791487762036
$ws.Range("K24").Value = 1.038423711281002
$ws.Range("L24").Value = 1.049341298043641
$ws.Range("M24").Value = 1.057140551695789
$ws.Range("N24").Value = 1.018738939840072
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.040676283627762
$ws.Range("D25").Value = 1.036024555526831
$ws.Range("E25").Value = 1.047895758235349
$ws.Range("F25").Value = 1.0560165592802
$ws.Range("I25").Value = 1.039419955054656
$ws.Range("J25").Value = 1.046239889502806
$ws.Range("K25").Value = 1.039076379263437
$ws.Range("L25").Value = 1.050910752972441
$ws.Range("M25").Value = 1.059006880281326
$ws.Range("N25").Value = 1.019230603469915
